$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 557.1
$ws.Cells.Item(5, 9).Value = 21.375
$ws.Cells.Item(5, 10).Value = 2700
$ws.Cells.Item(5, 11).Value = 21.375
$ws.Cells.Item(5, 12).Value = 2700
$ws.Cells.Item(5, 13).Value = 93.625
$ws.Cells.Item(5, 14).Value = -2930
$ws.Cells.Item(15, 8).Value = 89.72
$ws.Cells.Item(15, 9).Value = 89.72
$ws.Cells.Item(15, 11).Value = 269.16
$ws.Cells.Item(15, 13).Value = -100.16
$ws.Cells.Item(28, 8).Value = 322.29413
$ws.Cells.Item(28, 9).Value = 256.5
$ws.Cells.Item(28, 10).Value = 480.2
$ws.Cells.Item(28, 11).Value = 256.5
$ws.Cells.Item(28, 12).Value = 480.2
$ws.Cells.Item(28, 13).Value = 228.5
$ws.Cells.Item(28, 14).Value = -1450.2
$ws.Cells.Item(33, 8).Value = 183.91667
$ws.Cells.Item(33, 9).Value = 197
$ws.Cells.Item(33, 11).Value = 197
$ws.Cells.Item(33, 13).Value = 32
$ws.Cells.Item(100, 8).Value = 2161.7856
$ws.Cells.Item(100, 9).Value = 1253
$ws.Cells.Item(100, 10).Value = 2666.6667
$ws.Cells.Item(100, 11).Value = 1253
$ws.Cells.Item(100, 12).Value = 2666.6667
$ws.Cells.Item(100, 13).Value = -712
$ws.Cells.Item(100, 14).Value = -3748.6667
$ws.Cells.Item(129, 8).Value = 186188.11
$ws.Cells.Item(129, 10).Value = 218523.55
$ws.Cells.Item(129, 12).Value = 655570.6499999999
$ws.Cells.Item(129, 14).Value = -665570.6499999999
$ws.Cells.Item(138, 8).Value = 3399.973
$ws.Cells.Item(138, 9).Value = 5550
$ws.Cells.Item(138, 10).Value = 3139.3635
$ws.Cells.Item(138, 11).Value = 16650
$ws.Cells.Item(138, 12).Value = 9418.0905
$ws.Cells.Item(138, 13).Value = -11510
$ws.Cells.Item(138, 14).Value = -19698.0905
$ws.Cells.Item(141, 8).Value = 2877.4285
$ws.Cells.Item(141, 9).Value = 2334.875
$ws.Cells.Item(141, 10).Value = 3600.8333
$ws.Cells.Item(141, 11).Value = 7004.625
$ws.Cells.Item(141, 12).Value = 10802.4999
$ws.Cells.Item(141, 13).Value = -1824.625
$ws.Cells.Item(141, 14).Value = -21162.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1444.625
$ws.Cells.Item(2, 9).Value = 1470.1111
$ws.Cells.Item(2, 10).Value = 1411.8572
$ws.Cells.Item(2, 11).Value = 1470.1111
$ws.Cells.Item(2, 12).Value = 1411.8572
$ws.Cells.Item(2, 13).Value = -1357.1111
$ws.Cells.Item(2, 14).Value = -1637.8572
$ws.Cells.Item(32, 8).Value = 5757.35
$ws.Cells.Item(32, 9).Value = 4366.1094
$ws.Cells.Item(32, 11).Value = 4366.1094
$ws.Cells.Item(32, 13).Value = -4079.1094
$ws.Cells.Item(116, 8).Value = 1444.625
$ws.Cells.Item(116, 9).Value = 1470.1111
$ws.Cells.Item(116, 10).Value = 1411.8572
$ws.Cells.Item(116, 11).Value = 1470.1111
$ws.Cells.Item(116, 12).Value = 1411.8572
$ws.Cells.Item(116, 13).Value = 823.8888999999999
$ws.Cells.Item(116, 14).Value = -5999.8572
$ws.Cells.Item(122, 8).Value = 2021.7894
$ws.Cells.Item(122, 9).Value = 1845.2222
$ws.Cells.Item(122, 11).Value = 5535.6666
$ws.Cells.Item(122, 13).Value = -3085.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1444.625
$ws.Cells.Item(3, 9).Value = 1470.1111
$ws.Cells.Item(3, 10).Value = 1411.8572
$ws.Cells.Item(3, 11).Value = 1470.1111
$ws.Cells.Item(3, 12).Value = 1411.8572
$ws.Cells.Item(3, 13).Value = -1356.1111
$ws.Cells.Item(3, 14).Value = -1639.8572
$ws.Cells.Item(134, 8).Value = 6551.625
$ws.Cells.Item(134, 9).Value = 7272.5713
$ws.Cells.Item(134, 11).Value = 21817.7139
$ws.Cells.Item(134, 13).Value = -19282.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 58
$ws.Cells.Item(7, 10).Value = 48
$ws.Cells.Item(7, 12).Value = 48
$ws.Cells.Item(7, 14).Value = -274
$ws.Cells.Item(76, 8).Value = 25002250
$ws.Cells.Item(76, 9).Value = 25002250
$ws.Cells.Item(76, 11).Value = 25002250
$ws.Cells.Item(76, 13).Value = -25001935
$ws.Cells.Item(79, 8).Value = 25002250
$ws.Cells.Item(79, 9).Value = 25002250
$ws.Cells.Item(79, 11).Value = 25002250
$ws.Cells.Item(79, 13).Value = -25001158
$ws.Cells.Item(99, 8).Value = 3197.3928
$ws.Cells.Item(99, 9).Value = 2534.6191
$ws.Cells.Item(99, 11).Value = 2534.6191
$ws.Cells.Item(99, 13).Value = -1036.6191
$ws.Cells.Item(126, 8).Value = 3197.3928
$ws.Cells.Item(126, 9).Value = 2534.6191
$ws.Cells.Item(126, 11).Value = 7603.8573
$ws.Cells.Item(126, 13).Value = -5133.8573
$ws.Cells.Item(132, 8).Value = 5543.8335
$ws.Cells.Item(132, 9).Value = 4114.25
$ws.Cells.Item(132, 10).Value = 8403
$ws.Cells.Item(132, 11).Value = 12342.75
$ws.Cells.Item(132, 12).Value = 25209
$ws.Cells.Item(132, 13).Value = -9812.75
$ws.Cells.Item(132, 14).Value = -30269
$ws.Cells.Item(134, 8).Value = 1380
$ws.Cells.Item(134, 9).Value = 1133.3334
$ws.Cells.Item(134, 11).Value = 3400.0002
$ws.Cells.Item(134, 13).Value = -865.0001999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 3737.4
$ws.Cells.Item(3, 9).Value = 1647.5
$ws.Cells.Item(3, 10).Value = 12097
$ws.Cells.Item(3, 11).Value = 4942.5
$ws.Cells.Item(3, 12).Value = 36291
$ws.Cells.Item(3, 13).Value = -4830.5
$ws.Cells.Item(3, 14).Value = -36515
$ws.Cells.Item(5, 8).Value = 960.7037
$ws.Cells.Item(5, 9).Value = 573.19446
$ws.Cells.Item(5, 10).Value = 1735.7222
$ws.Cells.Item(5, 11).Value = 1719.58338
$ws.Cells.Item(5, 12).Value = 5207.1666
$ws.Cells.Item(5, 13).Value = -1607.58338
$ws.Cells.Item(5, 14).Value = -5431.1666
$ws.Cells.Item(103, 8).Value = 2293.2856
$ws.Cells.Item(103, 9).Value = 799.6
$ws.Cells.Item(103, 10).Value = 6027.5
$ws.Cells.Item(103, 11).Value = 2398.8
$ws.Cells.Item(103, 12).Value = 18082.5
$ws.Cells.Item(103, 13).Value = -1519.8
$ws.Cells.Item(103, 14).Value = -19840.5
$ws.Cells.Item(107, 8).Value = 3542.3333
$ws.Cells.Item(107, 10).Value = 254.09091
$ws.Cells.Item(107, 12).Value = 762.27273
$ws.Cells.Item(107, 14).Value = -4602.27273
$ws.Cells.Item(131, 8).Value = 722.8200000000001
$ws.Cells.Item(131, 10).Value = 728.9897
$ws.Cells.Item(131, 12).Value = 2186.9691
$ws.Cells.Item(131, 14).Value = -12266.9691
$ws.Cells.Item(135, 8).Value = 960.7037
$ws.Cells.Item(135, 9).Value = 573.19446
$ws.Cells.Item(135, 10).Value = 1735.7222
$ws.Cells.Item(135, 11).Value = 5158.75014
$ws.Cells.Item(135, 12).Value = 15621.4998
$ws.Cells.Item(135, 13).Value = -2623.75014
$ws.Cells.Item(135, 14).Value = -20691.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 55.210526
$ws.Cells.Item(2, 9).Value = 54.142857
$ws.Cells.Item(2, 10).Value = 58.2
$ws.Cells.Item(2, 11).Value = 54.142857
$ws.Cells.Item(2, 12).Value = 58.2
$ws.Cells.Item(2, 13).Value = 58.857143
$ws.Cells.Item(2, 14).Value = -284.2
$ws.Cells.Item(107, 8).Value = 683.2857
$ws.Cells.Item(107, 9).Value = 600
$ws.Cells.Item(107, 10).Value = 794.3333
$ws.Cells.Item(107, 11).Value = 600
$ws.Cells.Item(107, 12).Value = 794.3333
$ws.Cells.Item(107, 13).Value = 1320
$ws.Cells.Item(107, 14).Value = -4634.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(48, 8).Value = 14510.25
$ws.Cells.Item(48, 9).Value = 13041
$ws.Cells.Item(48, 10).Value = 15000
$ws.Cells.Item(48, 11).Value = 13041
$ws.Cells.Item(48, 12).Value = 15000
$ws.Cells.Item(48, 13).Value = -12380
$ws.Cells.Item(48, 14).Value = -16322
$ws.Cells.Item(93, 8).Value = 1208.3334
$ws.Cells.Item(93, 9).Value = 1051.2667
$ws.Cells.Item(93, 10).Value = 1601
$ws.Cells.Item(93, 11).Value = 1051.2667
$ws.Cells.Item(93, 12).Value = 1601
$ws.Cells.Item(93, 13).Value = 196.7333000000001
$ws.Cells.Item(93, 14).Value = -4097

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(105, 8).Value = 29600
$ws.Cells.Item(105, 10).Value = 29600
$ws.Cells.Item(105, 12).Value = 29600
$ws.Cells.Item(105, 14).Value = -36588
$ws.Cells.Item(107, 8).Value = 71428790
$ws.Cells.Item(107, 9).Value = 90909290
$ws.Cells.Item(107, 11).Value = 272727870
$ws.Cells.Item(107, 13).Value = -272725950
$ws.Cells.Item(122, 8).Value = 1427.0714
$ws.Cells.Item(122, 10).Value = 1553.2222
$ws.Cells.Item(122, 12).Value = 4659.6666
$ws.Cells.Item(122, 14).Value = -9559.6666
